$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Set the [Date] placeholder into B1 (next to the "Дата заполнения:" label in A1)
$ws.Range("B1").Value = "[Date]"

# Update the active selection to K8 (matches saved workbook state)
$ws.Range("K8").Select()
